$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the old thick-border blank
# separator row), shifting rows 2..6 down to 3..7 while keeping their
# values/styles/row-heights intact.
$ws.Rows.Item(2).Insert()

# ---- Fill in the new row 2 -------------------------------------------
# A2:C2 get the new "(percent)" captions (Kyrgyz / Russian / English).
$ws.Range("A2").Value = "(пайыз менен)"
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"

# Whole new row: no wrap, vertically centered, Times New Roman font.
$rowRange = $ws.Range("A2:K2")
$rowRange.WrapText = $false
$rowRange.VerticalAlignment = -4108
$rowRange.Font.Name = "Times New Roman"

# A2:C2: italic 8pt, horizontally + vertically centered caption style.
$rangeABC = $ws.Range("A2:C2")
$rangeABC.HorizontalAlignment = -4108
$rangeABC.Font.FontStyle = "Italic"
$rangeABC.Font.Size = 8

# Place the cursor/selection where the author left it after editing.
$ws.Range("C12").Select()
